$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 16.38160980024002
$ws.Cells.Item(2, 3).Value = 10.18239560898671
$ws.Cells.Item(2, 4).Value = 4.804074217130813
$ws.Cells.Item(2, 5).Value = 11.72429133267088
$ws.Cells.Item(2, 6).Value = 25.20946671743103
$ws.Cells.Item(2, 12).Value = 9.711364670693239
$ws.Cells.Item(2, 14).Value = 17.31939639268868
$ws.Cells.Item(2, 15).Value = 22.41874780881982

$ws.Cells.Item(3, 2).Value = 15.86753923791684
$ws.Cells.Item(3, 3).Value = 10.01151163393853
$ws.Cells.Item(3, 4).Value = 4.779721577372539
$ws.Cells.Item(3, 5).Value = 11.76251506379261
$ws.Cells.Item(3, 6).Value = 25.16341593570551
$ws.Cells.Item(3, 12).Value = 9.686041345469683
$ws.Cells.Item(3, 14).Value = 17.37509038927632
$ws.Cells.Item(3, 15).Value = 22.44492829154824

$ws.Cells.Item(4, 2).Value = 15.5458657298105
$ws.Cells.Item(4, 3).Value = 9.90426553881055
$ws.Cells.Item(4, 4).Value = 4.764547223297559
$ws.Cells.Item(4, 5).Value = 11.78801211125533
$ws.Cells.Item(4, 6).Value = 25.143107057851
$ws.Cells.Item(4, 12).Value = 9.672210771992139
$ws.Cells.Item(4, 14).Value = 17.41115727559874
$ws.Cells.Item(4, 15).Value = 22.4672224202852

$ws.Cells.Item(5, 2).Value = 15.41347028136338
$ws.Cells.Item(5, 3).Value = 9.860012208264775
$ws.Cells.Item(5, 4).Value = 4.758309742973792
$ws.Cells.Item(5, 5).Value = 11.79891187238898
$ws.Cells.Item(5, 6).Value = 25.13683911752402
$ws.Cells.Item(5, 12).Value = 9.667010711863437
$ws.Cells.Item(5, 14).Value = 17.42632624371685
$ws.Cells.Item(5, 15).Value = 22.47786733923504

$ws.Cells.Item(6, 2).Value = 15.3914132631475
$ws.Cells.Item(6, 3).Value = 9.852631805290203
$ws.Cells.Item(6, 4).Value = 4.757270821126591
$ws.Cells.Item(6, 5).Value = 11.80075253267929
$ws.Cells.Item(6, 6).Value = 25.13591971222768
$ws.Cells.Item(6, 12).Value = 9.666173688873959
$ws.Cells.Item(6, 14).Value = 17.42887354318452
$ws.Cells.Item(6, 15).Value = 22.47972900130826

$ws.Cells.Item(7, 2).Value = 15.54408522280159
$ws.Cells.Item(7, 3).Value = 9.903670902241812
$ws.Cells.Item(7, 4).Value = 4.764463317760629
$ws.Cells.Item(7, 5).Value = 11.78815704670752
$ws.Cells.Item(7, 6).Value = 25.14301439082782
$ws.Cells.Item(7, 12).Value = 9.672138871784121
$ws.Cells.Item(7, 14).Value = 17.41135993944956
$ws.Cells.Item(7, 15).Value = 22.46735967136167

$ws.Cells.Item(8, 2).Value = 16.20572898778028
$ws.Cells.Item(8, 3).Value = 10.12397711324735
$ws.Cells.Item(8, 4).Value = 4.795724615319326
$ws.Cells.Item(8, 5).Value = 11.73704986756083
$ws.Cells.Item(8, 6).Value = 25.19193855231814
$ws.Cells.Item(8, 12).Value = 9.702279116006553
$ws.Cells.Item(8, 14).Value = 17.33821208352207
$ws.Cells.Item(8, 15).Value = 22.4264818690059

$ws.Cells.Item(9, 2).Value = 17.44709633362995
$ws.Cells.Item(9, 3).Value = 10.53612478726997
$ws.Cells.Item(9, 4).Value = 4.855184124267863
$ws.Cells.Item(9, 5).Value = 11.65293086739459
$ws.Cells.Item(9, 6).Value = 25.35077301574046
$ws.Cells.Item(9, 12).Value = 9.774824328275713
$ws.Cells.Item(9, 14).Value = 17.20956346224825
$ws.Cells.Item(9, 15).Value = 22.39581291220275

$ws.Cells.Item(10, 2).Value = 18.31499493306902
$ws.Cells.Item(10, 3).Value = 10.82498361837261
$ws.Cells.Item(10, 4).Value = 4.897633674536299
$ws.Cells.Item(10, 5).Value = 11.60096446137827
$ws.Cells.Item(10, 6).Value = 25.50523847229332
$ws.Cells.Item(10, 12).Value = 9.83603700677849
$ws.Cells.Item(10, 14).Value = 17.12399764355785
$ws.Cells.Item(10, 15).Value = 22.40360266226525

$ws.Cells.Item(11, 2).Value = 18.69845852675738
$ws.Cells.Item(11, 3).Value = 10.95300663760591
$ws.Cells.Item(11, 4).Value = 4.916653289882637
$ws.Cells.Item(11, 5).Value = 11.57946346187436
$ws.Cells.Item(11, 6).Value = 25.58354111756756
$ws.Cells.Item(11, 12).Value = 9.865535496712411
$ws.Cells.Item(11, 14).Value = 17.08700162297917
$ws.Cells.Item(11, 15).Value = 22.41374512255915

$ws.Cells.Item(12, 2).Value = 18.8419061243457
$ws.Cells.Item(12, 3).Value = 11.0009716236232
$ws.Cells.Item(12, 4).Value = 4.923811672917263
$ws.Cells.Item(12, 5).Value = 11.57162949955315
$ws.Cells.Item(12, 6).Value = 25.61433008646883
$ws.Cells.Item(12, 12).Value = 9.876937203540944
$ws.Cells.Item(12, 14).Value = 17.07326848758353
$ws.Cells.Item(12, 15).Value = 22.41853435160063

$ws.Cells.Item(13, 2).Value = 18.81109239058177
$ws.Cells.Item(13, 3).Value = 10.99066482644747
$ws.Cells.Item(13, 4).Value = 4.922271974287917
$ws.Cells.Item(13, 5).Value = 11.5733029781708
$ws.Cells.Item(13, 6).Value = 25.60764885205515
$ws.Cells.Item(13, 12).Value = 9.87447146095522
$ws.Cells.Item(13, 14).Value = 17.07621388159733
$ws.Cells.Item(13, 15).Value = 22.41746073705292

$ws.Cells.Item(14, 2).Value = 18.71029604893467
$ws.Cells.Item(14, 3).Value = 10.95696321552785
$ws.Cells.Item(14, 4).Value = 4.91724310682964
$ws.Cells.Item(14, 5).Value = 11.57881278179713
$ws.Cells.Item(14, 6).Value = 25.58605145906719
$ws.Cells.Item(14, 12).Value = 9.866468924610917
$ws.Cells.Item(14, 14).Value = 17.0858662537808
$ws.Cells.Item(14, 15).Value = 22.41412013161662

$ws.Cells.Item(15, 2).Value = 18.64832233372133
$ws.Cells.Item(15, 3).Value = 10.93625217875815
$ws.Cells.Item(15, 4).Value = 4.914156996831938
$ws.Cells.Item(15, 5).Value = 11.58222781927851
$ws.Cells.Item(15, 6).Value = 25.57296999122017
$ws.Cells.Item(15, 12).Value = 9.861597063718474
$ws.Cells.Item(15, 14).Value = 17.09181458873178
$ws.Cells.Item(15, 15).Value = 22.41219740897689

$ws.Cells.Item(16, 2).Value = 18.2896949045513
$ws.Cells.Item(16, 3).Value = 10.81654652166827
$ws.Cells.Item(16, 4).Value = 4.896384664393759
$ws.Cells.Item(16, 5).Value = 11.60241267832618
$ws.Cells.Item(16, 6).Value = 25.50028141175169
$ws.Cells.Item(16, 12).Value = 9.834141927751967
$ws.Cells.Item(16, 14).Value = 17.12645415109041
$ws.Cells.Item(16, 15).Value = 22.40307260082283

$ws.Cells.Item(17, 2).Value = 18.06668512887175
$ws.Cells.Item(17, 3).Value = 10.74222442694634
$ws.Cells.Item(17, 4).Value = 4.885405909066028
$ws.Cells.Item(17, 5).Value = 11.61534354889982
$ws.Cells.Item(17, 6).Value = 25.45773525419227
$ws.Cells.Item(17, 12).Value = 9.817717921181316
$ws.Cells.Item(17, 14).Value = 17.14819768412915
$ws.Cells.Item(17, 15).Value = 22.39916495663608

$ws.Cells.Item(18, 2).Value = 17.93735452347722
$ws.Cells.Item(18, 3).Value = 10.69915975231866
$ws.Cells.Item(18, 4).Value = 4.879063963711205
$ws.Cells.Item(18, 5).Value = 11.62298232970384
$ws.Cells.Item(18, 6).Value = 25.43402138970391
$ws.Cells.Item(18, 12).Value = 9.808427213557735
$ws.Cells.Item(18, 14).Value = 17.1608855380464
$ws.Cells.Item(18, 15).Value = 22.39753858226135

$ws.Cells.Item(19, 2).Value = 17.89338732059304
$ws.Cells.Item(19, 3).Value = 10.68452531611363
$ws.Cells.Item(19, 4).Value = 4.876912073037501
$ws.Cells.Item(19, 5).Value = 11.62560324694334
$ws.Cells.Item(19, 6).Value = 25.4261229093143
$ws.Cells.Item(19, 12).Value = 9.805308512121519
$ws.Cells.Item(19, 14).Value = 17.16521263435205
$ws.Cells.Item(19, 15).Value = 22.39709460874864

$ws.Cells.Item(20, 2).Value = 18.09053572516746
$ws.Cells.Item(20, 3).Value = 10.7501691161663
$ws.Cells.Item(20, 4).Value = 4.886577452427536
$ws.Cells.Item(20, 5).Value = 11.61394620003767
$ws.Cells.Item(20, 6).Value = 25.46218608462889
$ws.Cells.Item(20, 12).Value = 9.819450189906849
$ws.Cells.Item(20, 14).Value = 17.14586426439381
$ws.Cells.Item(20, 15).Value = 22.39951664023446

$ws.Cells.Item(21, 2).Value = 18.73995115648635
$ws.Cells.Item(21, 3).Value = 10.96687637918014
$ws.Cells.Item(21, 4).Value = 4.918721416105727
$ws.Cells.Item(21, 5).Value = 11.57718605639836
$ws.Cells.Item(21, 6).Value = 25.59236442088508
$ws.Cells.Item(21, 12).Value = 9.86881323925223
$ws.Cells.Item(21, 14).Value = 17.08302362170176
$ws.Cells.Item(21, 15).Value = 22.41507561426766

$ws.Cells.Item(22, 2).Value = 19.15406241053286
$ws.Cells.Item(22, 3).Value = 11.10549584522914
$ws.Cells.Item(22, 4).Value = 4.93947231688106
$ws.Cells.Item(22, 5).Value = 11.55495663310814
$ws.Cells.Item(22, 6).Value = 25.68406362455992
$ws.Cells.Item(22, 12).Value = 9.902419684276548
$ws.Cells.Item(22, 14).Value = 17.04356465002658
$ws.Cells.Item(22, 15).Value = 22.43077203741115

$ws.Cells.Item(23, 2).Value = 18.93402697072784
$ws.Cells.Item(23, 3).Value = 11.03179653963262
$ws.Cells.Item(23, 4).Value = 4.928421380498945
$ws.Cells.Item(23, 5).Value = 11.56665646612094
$ws.Cells.Item(23, 6).Value = 25.63452279002153
$ws.Cells.Item(23, 12).Value = 9.884362423331353
$ws.Cells.Item(23, 14).Value = 17.06447750582902
$ws.Cells.Item(23, 15).Value = 22.42188914826833

$ws.Cells.Item(24, 2).Value = 18.07975634695653
$ws.Cells.Item(24, 3).Value = 10.746578366317
$ws.Cells.Item(24, 4).Value = 4.886047891297605
$ws.Cells.Item(24, 5).Value = 11.6145773040244
$ws.Cells.Item(24, 6).Value = 25.46017153755001
$ws.Cells.Item(24, 12).Value = 9.818666558297615
$ws.Cells.Item(24, 14).Value = 17.14691862023178
$ws.Cells.Item(24, 15).Value = 22.39935571219266

$ws.Cells.Item(25, 2).Value = 17.11836764928993
$ws.Cells.Item(25, 3).Value = 10.42695016197747
$ws.Cells.Item(25, 4).Value = 4.839307638973628
$ws.Cells.Item(25, 5).Value = 11.67396139963776
$ws.Cells.Item(25, 6).Value = 25.3011208861387
$ws.Cells.Item(25, 12).Value = 9.753787459348533
$ws.Cells.Item(25, 14).Value = 17.24278948507242
$ws.Cells.Item(25, 15).Value = 22.47972900130826
